$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (autogluon) - fill in previously empty results
$ws.Range("B3").Value = "0.341 (0.269 ± 0.053)"
$ws.Range("C3").Value = "00:01:42 (00:01:54 ± 00:00:08)"
$ws.Range("D3").Value = "00:00:00 (00:00:00 ± 00:00:00)"
$ws.Range("E3").Value = "[]"
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "5"
$ws.Range("F3").ClearFormats()

# Fix mojibake "Â±" -> "±" in rows 4, 6, 8 (columns B, C, D)
$ws.Range("B4").Value = "0.309 (0.278 ± 0.025)"
$ws.Range("C4").Value = "00:00:13 (00:00:17 ± 00:00:03)"
$ws.Range("D4").Value = "00:00:00 (00:00:00 ± 00:00:00)"

$ws.Range("B6").Value = "0.799 (0.716 ± 0.034)"
$ws.Range("C6").Value = "00:04:56 (00:05:02 ± 00:00:06)"
$ws.Range("D6").Value = "00:00:00 (00:00:00 ± 00:00:00)"

$ws.Range("B8").Value = "0.739 (0.661 ± 0.055)"
$ws.Range("C8").Value = "00:04:59 (00:05:00 ± 00:00:00)"
$ws.Range("D8").Value = "00:00:00 (00:00:00 ± 00:00:00)"
